# Fix quantities on 2WLS-30W-BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 100n ceramic capacitors qty 6 -> 4
$ws.Range("A2").Value = 4

# Row 4: 47u capacitors - parts list gains C9
$ws.Range("C4").Value = "C7, C8,C9"

# Row 5: BEAD inductors qty 3 -> 2
$ws.Range("A5").Value = 2

# Row 6: U1 5V DCDC Converter - drop "(see also Recom/Delta alternates on PCB)" suffix
$ws.Range("D6").Value = "5V DCDC Converter"

# Row 7: U2 15V DCDC Dual Converter qty 2 -> 1, drop suffix (keep trailing space)
$ws.Range("A7").Value = 1
$ws.Range("D7").Value = "15V DCDC Dual Converter "

# Update the active cell selection to match the saved view state
$ws.Range("F11").Select()
